$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">The example in Exhibit 6 shows a company with an initial higher short-term dividend</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">growth of</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:sSub>
          <m:e>
            <m:r>
              <m:t>g</m:t>
            </m:r>
          </m:e>
          <m:sub>
            <m:r>
              <m:t>s</m:t>
            </m:r>
          </m:sub>
        </m:sSub>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">for the first three years, followed by lower long-term growth</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:d>
          <m:dPr>
            <m:begChr m:val="("/>
            <m:sepChr m:val=""/>
            <m:endChr m:val=")"/>
            <m:grow/>
          </m:dPr>
          <m:e>
            <m:sSub>
              <m:e>
                <m:r>
                  <m:t>g</m:t>
                </m:r>
              </m:e>
              <m:sub>
                <m:r>
                  <m:t>l</m:t>
                </m:r>
              </m:sub>
            </m:sSub>
            <m:r>
              <m:rPr>
                <m:nor/>
                <m:sty m:val="p"/>
              </m:rPr>
              <m:t>, where </m:t>
            </m:r>
            <m:sSub>
              <m:e>
                <m:r>
                  <m:t>g</m:t>
                </m:r>
              </m:e>
              <m:sub>
                <m:r>
                  <m:t>s</m:t>
                </m:r>
              </m:sub>
            </m:sSub>
            <m:r>
              <m:rPr>
                <m:sty m:val="p"/>
              </m:rPr>
              <m:t>&gt;</m:t>
            </m:r>
            <m:sSub>
              <m:e>
                <m:r>
                  <m:t>g</m:t>
                </m:r>
              </m:e>
              <m:sub>
                <m:r>
                  <m:t>l</m:t>
                </m:r>
              </m:sub>
            </m:sSub>
          </m:e>
        </m:d>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">indefinitely thereafter. If we generalize the initial growth phase to</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:r>
          <m:t>n</m:t>
        </m:r>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">periods</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">followed by indefinite slower growth at</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:sSub>
          <m:e>
            <m:r>
              <m:t>g</m:t>
            </m:r>
          </m:e>
          <m:sub>
            <m:r>
              <m:t>l</m:t>
            </m:r>
          </m:sub>
        </m:sSub>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve">, we obtain a modified version of Equation</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">14 as follows:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <m:oMathPara>
        <m:oMathParaPr>
          <m:jc m:val="center"/>
        </m:oMathParaPr>
        <m:oMath>
          <m:r>
            <m:t>P</m:t>
          </m:r>
          <m:sSub>
            <m:e>
              <m:r>
                <m:t>V</m:t>
              </m:r>
            </m:e>
            <m:sub>
              <m:r>
                <m:t>t</m:t>
              </m:r>
            </m:sub>
          </m:sSub>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <m:t>=</m:t>
          </m:r>
          <m:nary>
            <m:naryPr>
              <m:chr m:val="∑"/>
              <m:limLoc m:val="subSup"/>
              <m:subHide m:val="off"/>
              <m:supHide m:val="off"/>
            </m:naryPr>
            <m:sub>
              <m:r>
                <m:t>i</m:t>
              </m:r>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>=</m:t>
              </m:r>
              <m:r>
                <m:t>1</m:t>
              </m:r>
            </m:sub>
            <m:sup>
              <m:r>
                <m:t>n</m:t>
              </m:r>
            </m:sup>
            <m:e>
              <m:f>
                <m:fPr>
                  <m:type m:val="bar"/>
                </m:fPr>
                <m:num>
                  <m:sSub>
                    <m:e>
                      <m:r>
                        <m:t>D</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:t>t</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:sSub>
                            <m:e>
                              <m:r>
                                <m:t>g</m:t>
                              </m:r>
                            </m:e>
                            <m:sub>
                              <m:r>
                                <m:t>s</m:t>
                              </m:r>
                            </m:sub>
                          </m:sSub>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>i</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:num>
                <m:den>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:r>
                            <m:t>r</m:t>
                          </m:r>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>i</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:den>
              </m:f>
            </m:e>
          </m:nary>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <m:t>+</m:t>
          </m:r>
          <m:nary>
            <m:naryPr>
              <m:chr m:val="∑"/>
              <m:limLoc m:val="subSup"/>
              <m:subHide m:val="off"/>
              <m:supHide m:val="off"/>
            </m:naryPr>
            <m:sub>
              <m:r>
                <m:t>j</m:t>
              </m:r>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>=</m:t>
              </m:r>
              <m:r>
                <m:t>n</m:t>
              </m:r>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>+</m:t>
              </m:r>
              <m:r>
                <m:t>1</m:t>
              </m:r>
            </m:sub>
            <m:sup>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>∞</m:t>
              </m:r>
            </m:sup>
            <m:e>
              <m:f>
                <m:fPr>
                  <m:type m:val="bar"/>
                </m:fPr>
                <m:num>
                  <m:sSub>
                    <m:e>
                      <m:r>
                        <m:t>D</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:t>t</m:t>
                      </m:r>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <m:t>+</m:t>
                      </m:r>
                      <m:r>
                        <m:t>n</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:sSub>
                            <m:e>
                              <m:r>
                                <m:t>g</m:t>
                              </m:r>
                            </m:e>
                            <m:sub>
                              <m:r>
                                <m:t>l</m:t>
                              </m:r>
                            </m:sub>
                          </m:sSub>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>j</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:num>
                <m:den>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:r>
                            <m:t>r</m:t>
                          </m:r>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>j</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:den>
              </m:f>
            </m:e>
          </m:nary>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="FirstParagraph"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">(15)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Note that the second expression in Equation 15 involves constant growth starting in</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:r>
          <m:t>n</m:t>
        </m:r>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">periods, for which we can substitute the geometric series simplification:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <m:oMathPara>
        <m:oMathParaPr>
          <m:jc m:val="center"/>
        </m:oMathParaPr>
        <m:oMath>
          <m:r>
            <m:t>P</m:t>
          </m:r>
          <m:sSub>
            <m:e>
              <m:r>
                <m:t>V</m:t>
              </m:r>
            </m:e>
            <m:sub>
              <m:r>
                <m:t>t</m:t>
              </m:r>
            </m:sub>
          </m:sSub>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <m:t>=</m:t>
          </m:r>
          <m:nary>
            <m:naryPr>
              <m:chr m:val="∑"/>
              <m:limLoc m:val="subSup"/>
              <m:subHide m:val="off"/>
              <m:supHide m:val="off"/>
            </m:naryPr>
            <m:sub>
              <m:r>
                <m:t>i</m:t>
              </m:r>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>=</m:t>
              </m:r>
              <m:r>
                <m:t>1</m:t>
              </m:r>
            </m:sub>
            <m:sup>
              <m:r>
                <m:t>n</m:t>
              </m:r>
            </m:sup>
            <m:e>
              <m:f>
                <m:fPr>
                  <m:type m:val="bar"/>
                </m:fPr>
                <m:num>
                  <m:sSub>
                    <m:e>
                      <m:r>
                        <m:t>D</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:t>t</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:sSub>
                            <m:e>
                              <m:r>
                                <m:t>g</m:t>
                              </m:r>
                            </m:e>
                            <m:sub>
                              <m:r>
                                <m:t>s</m:t>
                              </m:r>
                            </m:sub>
                          </m:sSub>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>i</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:num>
                <m:den>
                  <m:sSup>
                    <m:e>
                      <m:d>
                        <m:dPr>
                          <m:begChr m:val="("/>
                          <m:sepChr m:val=""/>
                          <m:endChr m:val=")"/>
                          <m:grow/>
                        </m:dPr>
                        <m:e>
                          <m:r>
                            <m:t>1</m:t>
                          </m:r>
                          <m:r>
                            <m:rPr>
                              <m:sty m:val="p"/>
                            </m:rPr>
                            <m:t>+</m:t>
                          </m:r>
                          <m:r>
                            <m:t>r</m:t>
                          </m:r>
                        </m:e>
                      </m:d>
                    </m:e>
                    <m:sup>
                      <m:r>
                        <m:t>i</m:t>
                      </m:r>
                    </m:sup>
                  </m:sSup>
                </m:den>
              </m:f>
            </m:e>
          </m:nary>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <m:t>+</m:t>
          </m:r>
          <m:f>
            <m:fPr>
              <m:type m:val="bar"/>
            </m:fPr>
            <m:num>
              <m:r>
                <m:t>E</m:t>
              </m:r>
              <m:d>
                <m:dPr>
                  <m:begChr m:val="("/>
                  <m:sepChr m:val=""/>
                  <m:endChr m:val=")"/>
                  <m:grow/>
                </m:dPr>
                <m:e>
                  <m:sSub>
                    <m:e>
                      <m:r>
                        <m:t>S</m:t>
                      </m:r>
                    </m:e>
                    <m:sub>
                      <m:r>
                        <m:t>t</m:t>
                      </m:r>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <m:t>+</m:t>
                      </m:r>
                      <m:r>
                        <m:t>n</m:t>
                      </m:r>
                    </m:sub>
                  </m:sSub>
                </m:e>
              </m:d>
            </m:num>
            <m:den>
              <m:sSup>
                <m:e>
                  <m:d>
                    <m:dPr>
                      <m:begChr m:val="("/>
                      <m:sepChr m:val=""/>
                      <m:endChr m:val=")"/>
                      <m:grow/>
                    </m:dPr>
                    <m:e>
                      <m:r>
                        <m:t>1</m:t>
                      </m:r>
                      <m:r>
                        <m:rPr>
                          <m:sty m:val="p"/>
                        </m:rPr>
                        <m:t>+</m:t>
                      </m:r>
                      <m:r>
                        <m:t>r</m:t>
                      </m:r>
                    </m:e>
                  </m:d>
                </m:e>
                <m:sup>
                  <m:r>
                    <m:t>n</m:t>
                  </m:r>
                </m:sup>
              </m:sSup>
            </m:den>
          </m:f>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="FirstParagraph"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">(16)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">where the stock value of the stock in n periods</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <m:oMath>
        <m:r>
          <m:rPr>
            <m:sty m:val="p"/>
          </m:rPr>
          <m:t>(</m:t>
        </m:r>
        <m:r>
          <m:t>E</m:t>
        </m:r>
        <m:d>
          <m:dPr>
            <m:begChr m:val="("/>
            <m:sepChr m:val=""/>
            <m:endChr m:val=")"/>
            <m:grow/>
          </m:dPr>
          <m:e>
            <m:sSub>
              <m:e>
                <m:r>
                  <m:t>S</m:t>
                </m:r>
              </m:e>
              <m:sub>
                <m:r>
                  <m:t>t</m:t>
                </m:r>
                <m:r>
                  <m:rPr>
                    <m:sty m:val="p"/>
                  </m:rPr>
                  <m:t>+</m:t>
                </m:r>
                <m:r>
                  <m:t>n</m:t>
                </m:r>
              </m:sub>
            </m:sSub>
          </m:e>
        </m:d>
      </m:oMath>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">is referred to as the terminal</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">value) and is equal to the following:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="BodyText"/>
      </w:pPr>
      <m:oMathPara>
        <m:oMathParaPr>
          <m:jc m:val="center"/>
        </m:oMathParaPr>
        <m:oMath>
          <m:r>
            <m:t>E</m:t>
          </m:r>
          <m:d>
            <m:dPr>
              <m:begChr m:val="("/>
              <m:sepChr m:val=""/>
              <m:endChr m:val=")"/>
              <m:grow/>
            </m:dPr>
            <m:e>
              <m:sSub>
                <m:e>
                  <m:r>
                    <m:t>S</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <m:t>t</m:t>
                  </m:r>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <m:t>+</m:t>
                  </m:r>
                  <m:r>
                    <m:t>n</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:e>
          </m:d>
          <m:r>
            <m:rPr>
              <m:sty m:val="p"/>
            </m:rPr>
            <m:t>=</m:t>
          </m:r>
          <m:f>
            <m:fPr>
              <m:type m:val="bar"/>
            </m:fPr>
            <m:num>
              <m:sSub>
                <m:e>
                  <m:r>
                    <m:t>D</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <m:t>t</m:t>
                  </m:r>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <m:t>+</m:t>
                  </m:r>
                  <m:r>
                    <m:t>n</m:t>
                  </m:r>
                  <m:r>
                    <m:rPr>
                      <m:sty m:val="p"/>
                    </m:rPr>
                    <m:t>+</m:t>
                  </m:r>
                  <m:r>
                    <m:t>1</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:num>
            <m:den>
              <m:r>
                <m:t>r</m:t>
              </m:r>
              <m:r>
                <m:rPr>
                  <m:sty m:val="p"/>
                </m:rPr>
                <m:t>−</m:t>
              </m:r>
              <m:sSub>
                <m:e>
                  <m:r>
                    <m:t>g</m:t>
                  </m:r>
                </m:e>
                <m:sub>
                  <m:r>
                    <m:t>l</m:t>
                  </m:r>
                </m:sub>
              </m:sSub>
            </m:den>
          </m:f>
        </m:oMath>
      </m:oMathPara>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="FirstParagraph"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">(17)</w:t>
      </w:r>
    </w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
